$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the FilesTab Neo4j query (B4): remove the File Type and Breed
# coalesce lines from the RETURN clause.
$newQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Irish Wolfhound'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
         coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newQuery

# The row height auto-adjusts in Excel because the cell text got shorter
# (two fewer wrapped lines in the query text). Set it explicitly to match.
$ws.Rows.Item(4).RowHeight = 217.5

# Update the selected cell in the sheet view from C4 to B4
$ws.Range("B4").Select()

$wb.Save()
